# Apply updated TPM-derived values to the NATMI LR-pairs sheet (Mdk-Lrp1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.180165333333334
$ws.Range("H2").Value = 6.540496
$ws.Range("I2").Value = 0.01970539991828544
$ws.Range("J2").Value = 0.01970539991828544
$ws.Range("M2").Value = 3.795192333333334
$ws.Range("N2").Value = 11.385577
$ws.Range("O2").Value = 0.01044213755712683
$ws.Range("P2").Value = 0.01044213755712683
$ws.Range("Q2").Value = 8.274146758465779
$ws.Range("R2").Value = 74.46732082619201
$ws.Range("S2").Value = 0.0002057664965649324
$ws.Range("T2").Value = 0.0002057664965649324

# Row 3
$ws.Range("G3").Value = 2.180165333333334
$ws.Range("H3").Value = 6.540496
$ws.Range("I3").Value = 0.01970539991828544
$ws.Range("J3").Value = 0.01970539991828544
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.6696287328350964
$ws.Range("P3").Value = 0.6696287328350964
$ws.Range("Q3").Value = 530.6007873245796
$ws.Range("R3").Value = 4775.407085921217
$ws.Range("S3").Value = 0.01319530197729029
$ws.Range("T3").Value = 0.01319530197729029

# Row 4
$ws.Range("G4").Value = 2.180165333333334
$ws.Range("H4").Value = 6.540496
$ws.Range("I4").Value = 0.01970539991828544
$ws.Range("J4").Value = 0.01970539991828544
$ws.Range("M4").Value = 29.801371
$ws.Range("N4").Value = 89.404113
$ws.Range("O4").Value = 0.08199584844219236
$ws.Range("P4").Value = 0.08199584844219235
$ws.Range("Q4").Value = 64.97191594000533
$ws.Range("R4").Value = 584.747243460048
$ws.Range("S4").Value = 0.001615760985192522
$ws.Range("T4").Value = 0.001615760985192522

# Row 5
$ws.Range("G5").Value = 2.180165333333334
$ws.Range("H5").Value = 6.540496
$ws.Range("I5").Value = 0.01970539991828544
$ws.Range("J5").Value = 0.01970539991828544
$ws.Range("M5").Value = 86.47679266666667
$ws.Range("N5").Value = 259.430378
$ws.Range("O5").Value = 0.2379332811655844
$ws.Range("P5").Value = 0.2379332811655844
$ws.Range("Q5").Value = 188.5337055097209
$ws.Range("R5").Value = 1696.803349587488
$ws.Range("S5").Value = 0.004688570459237692
$ws.Range("T5").Value = 0.004688570459237692

# Row 6
$ws.Range("I6").Value = 0.733713204346044
$ws.Range("J6").Value = 0.7337132043460441
$ws.Range("M6").Value = 3.795192333333334
$ws.Range("N6").Value = 11.385577
$ws.Range("O6").Value = 0.01044213755712683
$ws.Range("P6").Value = 0.01044213755712683
$ws.Range("Q6").Value = 308.0805645436292
$ws.Range("R6").Value = 2772.725080892662
$ws.Range("S6").Value = 0.007661534207261701
$ws.Range("T6").Value = 0.007661534207261702

# Row 7
$ws.Range("I7").Value = 0.733713204346044
$ws.Range("J7").Value = 0.7337132043460441
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.6696287328350964
$ws.Range("P7").Value = 0.6696287328350964
$ws.Range("Q7").Value = 19756.45282566409
$ws.Range("S7").Value = 0.4913154432906196
$ws.Range("T7").Value = 0.4913154432906197

# Row 8
$ws.Range("I8").Value = 0.733713204346044
$ws.Range("J8").Value = 0.7337132043460441
$ws.Range("M8").Value = 29.801371
$ws.Range("N8").Value = 89.404113
$ws.Range("O8").Value = 0.08199584844219236
$ws.Range("P8").Value = 0.08199584844219235
$ws.Range("Q8").Value = 2419.172045963275
$ws.Range("R8").Value = 21772.54841366948
$ws.Range("S8").Value = 0.06016143670359354
$ws.Range("T8").Value = 0.06016143670359354

# Row 9
$ws.Range("I9").Value = 0.733713204346044
$ws.Range("J9").Value = 0.7337132043460441
$ws.Range("M9").Value = 86.47679266666667
$ws.Range("N9").Value = 259.430378
$ws.Range("O9").Value = 0.2379332811655844
$ws.Range("P9").Value = 0.2379332811655844
$ws.Range("Q9").Value = 7019.886415419008
$ws.Range("R9").Value = 63178.97773877107
$ws.Range("S9").Value = 0.1745747901445692
$ws.Range("T9").Value = 0.1745747901445692

# Row 10
$ws.Range("G10").Value = 25.672264
$ws.Range("H10").Value = 77.016792
$ws.Range("I10").Value = 0.2320384702908474
$ws.Range("J10").Value = 0.2320384702908474
$ws.Range("M10").Value = 3.795192333333334
$ws.Range("N10").Value = 11.385577
$ws.Range("O10").Value = 0.01044213755712683
$ws.Range("P10").Value = 0.01044213755712683
$ws.Range("Q10").Value = 97.43117951210934
$ws.Range("R10").Value = 876.880615608984
$ws.Range("S10").Value = 0.002422977625322317
$ws.Range("T10").Value = 0.002422977625322317

# Row 11
$ws.Range("G11").Value = 25.672264
$ws.Range("H11").Value = 77.016792
$ws.Range("I11").Value = 0.2320384702908474
$ws.Range("J11").Value = 0.2320384702908474
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.6696287328350964
$ws.Range("P11").Value = 0.6696287328350964
$ws.Range("Q11").Value = 6248.023157939915
$ws.Range("R11").Value = 56232.20842145923
$ws.Range("S11").Value = 0.1553796268298543
$ws.Range("T11").Value = 0.1553796268298543

# Row 12
$ws.Range("G12").Value = 25.672264
$ws.Range("H12").Value = 77.016792
$ws.Range("I12").Value = 0.2320384702908474
$ws.Range("J12").Value = 0.2320384702908474
$ws.Range("M12").Value = 29.801371
$ws.Range("N12").Value = 89.404113
$ws.Range("O12").Value = 0.08199584844219236
$ws.Range("P12").Value = 0.08199584844219235
$ws.Range("Q12").Value = 765.0686638739439
$ws.Range("R12").Value = 6885.617974865495
$ws.Range("S12").Value = 0.01902619124272648
$ws.Range("T12").Value = 0.01902619124272648

# Row 13
$ws.Range("G13").Value = 25.672264
$ws.Range("H13").Value = 77.016792
$ws.Range("I13").Value = 0.2320384702908474
$ws.Range("J13").Value = 0.2320384702908474
$ws.Range("M13").Value = 86.47679266666667
$ws.Range("N13").Value = 259.430378
$ws.Range("O13").Value = 0.2379332811655844
$ws.Range("P13").Value = 0.2379332811655844
$ws.Range("Q13").Value = 2220.055051211931
$ws.Range("R13").Value = 19980.49546090738
$ws.Range("S13").Value = 0.05520967459294429
$ws.Range("T13").Value = 0.05520967459294429

# Row 14
$ws.Range("G14").Value = 1.608999666666667
$ws.Range("H14").Value = 4.826999
$ws.Range("I14").Value = 0.01454292544482312
$ws.Range("J14").Value = 0.01454292544482312
$ws.Range("M14").Value = 3.795192333333334
$ws.Range("N14").Value = 11.385577
$ws.Range("O14").Value = 0.01044213755712683
$ws.Range("P14").Value = 0.01044213755712683
$ws.Range("Q14").Value = 6.106463199269223
$ws.Range("R14").Value = 54.95816879342301
$ws.Range("S14").Value = 0.000151859227977883
$ws.Range("T14").Value = 0.000151859227977883

# Row 15
$ws.Range("G15").Value = 1.608999666666667
$ws.Range("H15").Value = 4.826999
$ws.Range("I15").Value = 0.01454292544482312
$ws.Range("J15").Value = 0.01454292544482312
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.6696287328350964
$ws.Range("P15").Value = 0.6696287328350964
$ws.Range("Q15").Value = 391.5925443292005
$ws.Range("R15").Value = 3524.332898962804
$ws.Range("S15").Value = 0.00973836073733219
$ws.Range("T15").Value = 0.00973836073733219

# Row 16
$ws.Range("G16").Value = 1.608999666666667
$ws.Range("H16").Value = 4.826999
$ws.Range("I16").Value = 0.01454292544482312
$ws.Range("J16").Value = 0.01454292544482312
$ws.Range("M16").Value = 29.801371
$ws.Range("N16").Value = 89.404113
$ws.Range("O16").Value = 0.08199584844219236
$ws.Range("P16").Value = 0.08199584844219235
$ws.Range("Q16").Value = 47.95039600520966
$ws.Range("R16").Value = 431.5535640468869
$ws.Range("S16").Value = 0.00119245951067982
$ws.Range("T16").Value = 0.00119245951067982

# Row 17
$ws.Range("G17").Value = 1.608999666666667
$ws.Range("H17").Value = 4.826999
$ws.Range("I17").Value = 0.01454292544482312
$ws.Range("J17").Value = 0.01454292544482312
$ws.Range("M17").Value = 86.47679266666667
$ws.Range("N17").Value = 259.430378
$ws.Range("O17").Value = 0.2379332811655844
$ws.Range("P17").Value = 0.2379332811655844
$ws.Range("Q17").Value = 139.1411305750691
$ws.Range("R17").Value = 1252.270175175622
$ws.Range("S17").Value = 0.003460245968833232
$ws.Range("T17").Value = 0.003460245968833232

